$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46060 -> 46061) for every data row (rows 2 through 412).
$firstRow = 2
$lastRow = 412

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 46060) {
        $cell.Value = 46061
    }
}
